$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Starting Era" values (B2:B15) from 340 to 339 ---
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Value = 339
}

# --- Add the new V41 models (rows 16-23), only considering models since
#     the start of the daily tournament ---
$newModels = @(
    "LG_LGBM_V41_CYRUS20",
    "LG_LGBM_V41_CAROLI20",
    "LG_LGBM_V41_XERXES20",
    "LG_LGBM_V41_SAM20",
    "LG_LGBM_V41_CYRUS60",
    "LG_LGBM_V41_CAROLI60",
    "LG_LGBM_V41_XERXES60",
    "LG_LGBM_V41_SAM60"
)

$row = 16
foreach ($name in $newModels) {
    $ws.Cells.Item($row, 1).Value = $name
    # Match the font/style used by the other highlighted ModelName cells
    # (A6:A9, A12:A15) so the new rows get the same cell style (s="1").
    $ws.Cells.Item($row, 1).Font.Color = $ws.Cells.Item(6, 1).Font.Color

    $ws.Cells.Item($row, 2).Value = 339
    $ws.Cells.Item($row, 3).Value = 1000
    $ws.Cells.Item($row, 4).Value = "-"

    $row++
}

# --- Update the view: scroll so row 3 is at top, and select F19 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("F19").Select()
